$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.536023666666666
$ws.Range("H2").Value = 7.608070999999999
$ws.Range("I2").Value = 0.02556389501525096
$ws.Range("J2").Value = 0.02556389501525096
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 124.0518796666667
$ws.Range("N2").Value = 372.155639
$ws.Range("O2").Value = 0.2841471011719914
$ws.Range("P2").Value = 0.2841471011719915
$ws.Range("Q2").Value = 314.598502729152
$ws.Range("R2").Value = 2831.386524562368
$ws.Range("S2").Value = 0.007263906663248682
$ws.Range("T2").Value = 0.007263906663248685
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.536023666666666
$ws.Range("H3").Value = 7.608070999999999
$ws.Range("I3").Value = 0.02556389501525096
$ws.Range("J3").Value = 0.02556389501525096
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 216.8516133333333
$ws.Range("N3").Value = 650.55484
$ws.Range("O3").Value = 0.4967095821418112
$ws.Range("P3").Value = 0.4967095821418113
$ws.Range("Q3").Value = 549.9408235681822
$ws.Range("R3").Value = 4949.46741211364
$ws.Range("S3").Value = 0.01269783161094244
$ws.Range("T3").Value = 0.01269783161094244
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.536023666666666
$ws.Range("H4").Value = 7.608070999999999
$ws.Range("I4").Value = 0.02556389501525096
$ws.Range("J4").Value = 0.02556389501525096
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 95.67277033333333
$ws.Range("N4").Value = 287.018311
$ws.Range("O4").Value = 0.2191433166861973
$ws.Range("P4").Value = 0.2191433166861973
$ws.Range("Q4").Value = 242.6284098208979
$ws.Range("R4").Value = 2183.655688388081
$ws.Range("S4").Value = 0.005602156741059841
$ws.Range("T4").Value = 0.005602156741059842
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 75.11538433333334
$ws.Range("H5").Value = 225.346153
$ws.Range("I5").Value = 0.7571860721834327
$ws.Range("J5").Value = 0.7571860721834328
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 124.0518796666667
$ws.Range("N5").Value = 372.155639
$ws.Range("O5").Value = 0.2841471011719914
$ws.Range("P5").Value = 0.2841471011719915
$ws.Range("Q5").Value = 9318.204618434085
$ws.Range("R5").Value = 83863.84156590677
$ws.Range("S5").Value = 0.2151522274587287
$ws.Range("T5").Value = 0.2151522274587287
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 75.11538433333334
$ws.Range("H6").Value = 225.346153
$ws.Range("I6").Value = 0.7571860721834327
$ws.Range("J6").Value = 0.7571860721834328
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 216.8516133333333
$ws.Range("N6").Value = 650.55484
$ws.Range("O6").Value = 0.4967095821418112
$ws.Range("P6").Value = 0.4967095821418113
$ws.Range("Q6").Value = 16288.89227883673
$ws.Range("R6").Value = 146600.0305095305
$ws.Range("S6").Value = 0.3761015775178322
$ws.Range("T6").Value = 0.3761015775178322
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 75.11538433333334
$ws.Range("H7").Value = 225.346153
$ws.Range("I7").Value = 0.7571860721834327
$ws.Range("J7").Value = 0.7571860721834328
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 95.67277033333333
$ws.Range("N7").Value = 287.018311
$ws.Range("O7").Value = 0.2191433166861973
$ws.Range("P7").Value = 0.2191433166861973
$ws.Range("Q7").Value = 7186.496913823065
$ws.Range("R7").Value = 64678.47222440758
$ws.Range("S7").Value = 0.1659322672068718
$ws.Range("T7").Value = 0.1659322672068719
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.551928
$ws.Range("H8").Value = 64.655784
$ws.Range("I8").Value = 0.2172500328013163
$ws.Range("J8").Value = 0.2172500328013163
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 124.0518796666667
$ws.Range("N8").Value = 372.155639
$ws.Range("O8").Value = 0.2841471011719914
$ws.Range("P8").Value = 0.2841471011719915
$ws.Range("Q8").Value = 2673.557178840664
$ws.Range("R8").Value = 24062.01460956597
$ws.Range("S8").Value = 0.06173096705001407
$ws.Range("T8").Value = 0.06173096705001409
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.551928
$ws.Range("H9").Value = 64.655784
$ws.Range("I9").Value = 0.2172500328013163
$ws.Range("J9").Value = 0.2172500328013163
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 216.8516133333333
$ws.Range("N9").Value = 650.55484
$ws.Range("O9").Value = 0.4967095821418112
$ws.Range("P9").Value = 0.4967095821418113
$ws.Range("Q9").Value = 4673.570357243841
$ws.Range("R9").Value = 42062.13321519456
$ws.Range("S9").Value = 0.1079101730130366
$ws.Range("T9").Value = 0.1079101730130366
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 21.551928
$ws.Range("H10").Value = 64.655784
$ws.Range("I10").Value = 0.2172500328013163
$ws.Range("J10").Value = 0.2172500328013163
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 95.67277033333333
$ws.Range("N10").Value = 287.018311
$ws.Range("O10").Value = 0.2191433166861973
$ws.Range("P10").Value = 0.2191433166861973
$ws.Range("Q10").Value = 2061.932657784536
$ws.Range("R10").Value = 18557.39392006082
$ws.Range("S10").Value = 0.04760889273826559
$ws.Range("T10").Value = 0.0476088927382656
